# Generate Report for Handback
# Update the timestamp strings recorded on the handback-status report.
$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the synced file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 23:03:21"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 23:03:16"
$wsZhCn.Range("K2").Value = "2016-08-22 23:03:41"

# "de-de" sheet: Correspond Handback DateTime (Correspond Handoff Datetime,
# which mirrors the Overview value, is updated via the shared value above).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 23:03:21"
$wsDeDe.Range("K2").Value = "2016-08-22 23:03:48"
